$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5408568.5
$ws.Range("I137").Value = 1694.3846
$ws.Range("J137").Value = 18188452
$ws.Range("K137").Value = 5083.1538
$ws.Range("L137").Value = 54565356
$ws.Range("M137").Value = -2533.1538
$ws.Range("N137").Value = -54570456
$ws.Range("H138").Value = 2843112.5
$ws.Range("I138").Value = 2134.9333
$ws.Range("J138").Value = 4312583.5
$ws.Range("K138").Value = 6404.7999
$ws.Range("L138").Value = 12937750.5
$ws.Range("M138").Value = -1264.7999
$ws.Range("N138").Value = -12948030.5
$ws.Range("H141").Value = 3724.75
$ws.Range("I141").Value = 3466.3333
$ws.Range("J141").Value = 4500
$ws.Range("K141").Value = 10398.9999
$ws.Range("L141").Value = 13500
$ws.Range("M141").Value = -5218.999899999999
$ws.Range("N141").Value = -23860

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3136.5
$ws.Range("I45").Value = 3417.7144
$ws.Range("J45").Value = 2742.8
$ws.Range("K45").Value = 3417.7144
$ws.Range("L45").Value = 2742.8
$ws.Range("M45").Value = -3040.7144
$ws.Range("N45").Value = -3496.8
$ws.Range("H97").Value = 5455.273
$ws.Range("I97").Value = 7615.8667
$ws.Range("J97").Value = 825.4286
$ws.Range("K97").Value = 7615.8667
$ws.Range("L97").Value = 825.4286
$ws.Range("M97").Value = -7119.8667
$ws.Range("N97").Value = -1817.4286
$ws.Range("H122").Value = 5991.609
$ws.Range("I122").Value = 7405.5293
$ws.Range("J122").Value = 1985.5
$ws.Range("K122").Value = 22216.5879
$ws.Range("L122").Value = 5956.5
$ws.Range("M122").Value = -19766.5879
$ws.Range("N122").Value = -10856.5
$ws.Range("H132").Value = 8067106.5
$ws.Range("I132").Value = 17859356
$ws.Range("J132").Value = 2901.2354
$ws.Range("K132").Value = 53578068
$ws.Range("L132").Value = 8703.706200000001
$ws.Range("M132").Value = -53575538
$ws.Range("N132").Value = -13763.7062

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 7993.923
$ws.Range("I80").Value = 15630
$ws.Range("J80").Value = 357.84616
$ws.Range("K80").Value = 15630
$ws.Range("L80").Value = 357.84616
$ws.Range("M80").Value = -14632
$ws.Range("N80").Value = -2353.84616
$ws.Range("H83").Value = 7993.923
$ws.Range("I83").Value = 15630
$ws.Range("J83").Value = 357.84616
$ws.Range("K83").Value = 78150
$ws.Range("L83").Value = 1789.2308
$ws.Range("M83").Value = -73158
$ws.Range("N83").Value = -11773.2308
$ws.Range("H134").Value = 4517.4766
$ws.Range("I134").Value = 3687.0317
$ws.Range("J134").Value = 6792.174
$ws.Range("K134").Value = 11061.0951
$ws.Range("L134").Value = 20376.522
$ws.Range("M134").Value = -8526.0951
$ws.Range("N134").Value = -25446.522

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2334.7188
$ws.Range("I58").Value = 1021.5
$ws.Range("J58").Value = 4523.4165
$ws.Range("K58").Value = 1021.5
$ws.Range("L58").Value = 4523.4165
$ws.Range("M58").Value = -818.5
$ws.Range("N58").Value = -4929.4165
$ws.Range("H132").Value = 1762.2245
$ws.Range("I132").Value = 1489.1316
$ws.Range("J132").Value = 2705.6365
$ws.Range("K132").Value = 4467.3948
$ws.Range("L132").Value = 8116.9095
$ws.Range("M132").Value = -1937.3948
$ws.Range("N132").Value = -13176.9095
$ws.Range("H136").Value = 2334.7188
$ws.Range("I136").Value = 1021.5
$ws.Range("J136").Value = 4523.4165
$ws.Range("K136").Value = 3064.5
$ws.Range("L136").Value = 13570.2495
$ws.Range("M136").Value = -514.5
$ws.Range("N136").Value = -18670.2495

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1172.2644
$ws.Range("I68").Value = 874.4783
$ws.Range("J68").Value = 1506.3658
$ws.Range("K68").Value = 2623.4349
$ws.Range("L68").Value = 4519.097400000001
$ws.Range("M68").Value = -1812.4349
$ws.Range("N68").Value = -6141.097400000001
$ws.Range("H71").Value = 1172.2644
$ws.Range("I71").Value = 874.4783
$ws.Range("J71").Value = 1506.3658
$ws.Range("K71").Value = 7870.3047
$ws.Range("L71").Value = 13557.2922
$ws.Range("M71").Value = -3814.3047
$ws.Range("N71").Value = -21669.2922
$ws.Range("H125").Value = 4655.5
$ws.Range("J125").Value = 4655.5
$ws.Range("L125").Value = 13966.5
$ws.Range("N125").Value = -23806.5
$ws.Range("H131").Value = 864.92
$ws.Range("I131").Value = 278.33334
$ws.Range("J131").Value = 902.3617
$ws.Range("K131").Value = 835.0000200000001
$ws.Range("L131").Value = 2707.0851
$ws.Range("M131").Value = 4204.99998
$ws.Range("N131").Value = -12787.0851
$ws.Range("H134").Value = 4761.579
$ws.Range("I134").Value = 3557
$ws.Range("J134").Value = 6100
$ws.Range("K134").Value = 10671
$ws.Range("L134").Value = 18300
$ws.Range("M134").Value = -5601
$ws.Range("N134").Value = -28440
$ws.Range("H139").Value = 2469.3333
$ws.Range("I139").Value = 1844
$ws.Range("J139").Value = 2782
$ws.Range("K139").Value = 5532
$ws.Range("L139").Value = 8346
$ws.Range("M139").Value = -392
$ws.Range("N139").Value = -18626

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1830.6923
$ws.Range("I61").Value = 1755.5555
$ws.Range("J61").Value = 1999.75
$ws.Range("K61").Value = 1755.5555
$ws.Range("L61").Value = 1999.75
$ws.Range("M61").Value = -1553.5555
$ws.Range("N61").Value = -2403.75
$ws.Range("H113").Value = 1830.6923
$ws.Range("I113").Value = 1755.5555
$ws.Range("J113").Value = 1999.75
$ws.Range("K113").Value = 1755.5555
$ws.Range("L113").Value = 1999.75
$ws.Range("M113").Value = 414.4445000000001
$ws.Range("N113").Value = -6339.75

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2245.9546
$ws.Range("I113").Value = 2402.2
$ws.Range("J113").Value = 2115.75
$ws.Range("K113").Value = 7206.599999999999
$ws.Range("L113").Value = 6347.25
$ws.Range("M113").Value = -5036.599999999999
$ws.Range("N113").Value = -10687.25
$ws.Range("H132").Value = 4778.6387
$ws.Range("I132").Value = 5686.08
$ws.Range("K132").Value = 17058.24
$ws.Range("M132").Value = -14528.24
$ws.Range("H136").Value = 1229.5
$ws.Range("I136").Value = 1145.4286
$ws.Range("J136").Value = 2995
$ws.Range("K136").Value = 3436.2858
$ws.Range("L136").Value = 8985
$ws.Range("M136").Value = -886.2857999999997
$ws.Range("N136").Value = -14085
